$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.379.02"
$ws.Range("E2").Value = "  +0.38%  "

# Row 3
$ws.Range("D3").Value = "1.871.17"
$ws.Range("E3").Value = "  -0.70%  "

# Row 4
$ws.Range("E4").Value = "  +0.87%  "

# Row 5
$ws.Range("D5").Value = "'316.57"
$ws.Range("E5").Value = "  +0.30%  "

# Row 6
$ws.Range("D6").Value = "'1.017"
$ws.Range("E6").Value = "  +0.71%  "

# Row 7
$ws.Range("D7").Value = "'0.5108"
$ws.Range("E7").Value = "  -0.66%  "

# Row 8
$ws.Range("D8").Value = "'0.3959"
$ws.Range("E8").Value = "  +1.31%  "

# Row 9
$ws.Range("D9").Value = "'0.08493"
$ws.Range("E9").Value = "  +1.13%  "

# Row 10
$ws.Range("D10").Value = "'1.110"
$ws.Range("E10").Value = "  -1.23%  "

# Row 11
$ws.Range("D11").Value = "'6.248"
$ws.Range("E11").Value = "  -0.30%  "

# Row 12
$ws.Range("D12").Value = "'20.46"
$ws.Range("E12").Value = "  -0.91%  "

# Row 13
$ws.Range("D13").Value = "1.825.50"
$ws.Range("E13").Value = "  -2.77%  "

# Row 14
$ws.Range("B14").Value = "BinanceUSD"
$ws.Range("C14").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D14").Value = "'1.021"
$ws.Range("E14").Value = "  +1.40%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.221"
$ws.Range("E15").Value = "  -0.93%  "

# Row 16
$ws.Range("D16").Value = "'0.00001113"
$ws.Range("E16").Value = "  +0.61%  "

# Row 17
$ws.Range("D17").Value = "'90.57"
$ws.Range("E17").Value = "  -0.63%  "

# Row 18
$ws.Range("D18").Value = "'0.06770"
$ws.Range("E18").Value = "  +0.85%  "

# Row 19
$ws.Range("D19").Value = "'17.68"

# Row 20
$ws.Range("D20").Value = "'1.016"
$ws.Range("E20").Value = "  +0.72%  "

# Row 21
$ws.Range("D21").Value = "'5.949"
$ws.Range("E21").Value = "  -1.62%  "

# Row 22
$ws.Range("D22").Value = "28.390.99"

# Row 23
$ws.Range("D23").Value = "'11.14"
$ws.Range("E23").Value = "  -0.27%  "

# Row 24
$ws.Range("D24").Value = "'2.284"
$ws.Range("E24").Value = "  +0.74%  "

# Row 25
$ws.Range("D25").Value = "'162.14"
$ws.Range("E25").Value = "  +1.48%  "

# Row 26
$ws.Range("D26").Value = "2.030.02"
$ws.Range("E26").Value = "  -2.80%  "

# Row 27
$ws.Range("D27").Value = "'20.76"
$ws.Range("E27").Value = "  +0.25%  "

# Row 28
$ws.Range("D28").Value = "'2.364"
$ws.Range("E28").Value = "  -4.16%  "

# Row 30
$ws.Range("D30").Value = "'0.1051"
$ws.Range("E30").Value = "  -0.69%  "

# Row 31
$ws.Range("D31").Value = "'1.036"
$ws.Range("E31").Value = "  -0.28%  "

# Row 32
$ws.Range("D32").Value = "'5.764"
$ws.Range("E32").Value = "  -2.02%  "

# Row 33
$ws.Range("D33").Value = "'3.644"
$ws.Range("E33").Value = "  +0.47%  "

# Row 34
$ws.Range("D34").Value = "'0.02423"
$ws.Range("E34").Value = "  -1.17%  "

# Row 35
$ws.Range("D35").Value = "'0.06455"
$ws.Range("E35").Value = "  -1.93%  "

# Row 36
$ws.Range("D36").Value = "'0.2188"
$ws.Range("E36").Value = "  -1.44%  "

# Row 37
$ws.Range("D37").Value = "'8.853"
$ws.Range("E37").Value = "  -7.59%  "

# Row 38
$ws.Range("D38").Value = "'1.259"
$ws.Range("E38").Value = "  +1.10%  "

# Row 39
$ws.Range("D39").Value = "'1.180"
$ws.Range("E39").Value = "  -1.62%  "

# Row 40
$ws.Range("D40").Value = "'0.6363"
$ws.Range("E40").Value = "  -2.06%  "

# Row 41
$ws.Range("D41").Value = "'4.980"

# Row 42
$ws.Range("D42").Value = "'11.24"
$ws.Range("E42").Value = "  -0.83%  "

# Row 43
$ws.Range("D43").Value = "'0.6023"
$ws.Range("E43").Value = "  -1.50%  "

# Row 44
$ws.Range("D44").Value = "'13.02"
$ws.Range("E44").Value = "  -0.81%  "

# Row 45
$ws.Range("D45").Value = "'3.699"
$ws.Range("E45").Value = "  -0.05%  "

# Row 46
$ws.Range("D46").Value = "'1.214"
$ws.Range("E46").Value = "  -5.23%  "

# Row 47
$ws.Range("D47").Value = "'1.991"
$ws.Range("E47").Value = "  -1.41%  "

# Row 48
$ws.Range("D48").Value = "'1.201"
$ws.Range("E48").Value = "  -2.85%  "

# Row 49
$ws.Range("D49").Value = "'120.82"
$ws.Range("E49").Value = "  -0.30%  "

# Row 50
$ws.Range("D50").Value = "'0.06843"
$ws.Range("E50").Value = "  -1.30%  "

# Row 51
$ws.Range("D51").Value = "'76.29"
$ws.Range("E51").Value = "  -2.27%  "
